$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen the "Key" column (C) so the new, longer keys are readable ---
# (column C grows from 32.625 to ~61.875 "characters"; COM only exposes the
# quantized ColumnWidth setter, so this lands on the closest representable
# pixel-quantized width)
$ws.Columns.Item(3).ColumnWidth = 61.14

# --- Row 24: Exception_get_summary_with_reason -------------------------
# Clone formatting from the last existing data row (23) so styles/number
# formats match, then overwrite with the new row's values.
$ws.Range("A23:C23").Copy()
$ws.Range("A24:C24").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E23:H23").Copy()
$ws.Range("E24:H24").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A24").Value = "CoreLib, Private.CoreLib"
$ws.Range("B24").Value = "Strings"
$ws.Range("C24").Value = "Exception_get_summary_with_reason"
$ws.Range("E24").Value = 'When "{0}", an exception of type "{1}" is thrown: {2}'
$ws.Range("F24").Value = 'When "{0}", an exception of type "{1}" is thrown: {2}'
$ws.Range("G24").Value = '當 「{0}」 時，引發了一個 「{1}」 類型的異常：{2}'
$ws.Range("H24").Value = '当 “{0}” 时，引发了一个 “{1}” 类型的异常：{2}'

# --- Row 25: Exception_get_summary_without_reason -----------------------
$ws.Range("A23:C23").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E23:H23").Copy()
$ws.Range("E25:H25").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A25").Value = "CoreLib, Private.CoreLib"
$ws.Range("B25").Value = "Strings"
$ws.Range("C25").Value = "Exception_get_summary_without_reason"
$ws.Range("E25").Value = 'When method "{0}" is called, an exception of type "{1}" is thrown: {2}'
$ws.Range("F25").Value = 'When method "{0}" is called, an exception of type "{1}" is thrown: {2}'
$ws.Range("G25").Value = '當調用方法 「{0}」 時，引發了一個 「{1}」 類型的異常：{2}'
$ws.Range("H25").Value = '当调用方法 “{0}” 时，引发了一个 “{1}” 类型的异常：{2}'

# --- Update the selection to match the new authored state ---------------
$ws.Range("A23:B25").Select() | Out-Null
